# Update the "想去人数" (interested-count) figures in column F for both the
# "展览" and "全部类型" worksheets. Both sheets mirror the same rows, so the
# same set of row/value updates is applied to each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 12556
    13 = 13648
    14 = 13958
    19 = 1049
    22 = 499
    23 = 5060
    24 = 249
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
